# Adds a visualization-friendly header row (Title Case / spaced-out column
# names) to both sheets of the "active features" report, replacing the raw
# snake_case field names that used to live in A1:E1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$headers = @("Layer", "Lang", "Feature ID", "Interpretation", "Tokens")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'A' + $i)
    $ws1.Range("$col`1").Value = $headers[$i]
    $ws2.Range("$col`1").Value = $headers[$i]
}

# Restore the view/selection state: sheet1's selection is set first, then
# sheet2's, so sheet2 (the second tab) ends up as the active/selected sheet,
# matching the saved workbook view.
$ws1.Range("J2").Select() | Out-Null
$ws2.Range("G2").Select() | Out-Null
